$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match header formatting (bold, centered, top-aligned, thin border) used by B1:H1
$hdr = $ws.Range("I1:J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# Data values for I and J columns, rows 2-66
$values = @{
    2 = @(9, 9)
    3 = @(9, 9)
    4 = @(8, 8)
    5 = @(8, 8)
    6 = @(7, 7)
    7 = @(9, 9)
    8 = @(8, 8)
    9 = @(8, 8)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(7, 8)
    16 = @(8, 8)
    17 = @(9, 9)
    18 = @(8, 8)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(8, 8)
    24 = @(8, 8)
    25 = @(8, 8)
    26 = @(9, 9)
    27 = @(8, 8)
    28 = @(7, 7)
    29 = @(8, 8)
    30 = @(9, 9)
    31 = @(6, 7)
    32 = @(6, 7)
    33 = @(7, 7)
    34 = @(8, 8)
    35 = @(7, 7)
    36 = @(5, 6)
    37 = @(3, 4)
    38 = @(8, 8)
    39 = @(10, 10)
    40 = @(10, 10)
    41 = @(7, 8)
    42 = @(7, 8)
    43 = @(6, 7)
    44 = @(6, 7)
    45 = @(9, 9)
    46 = @(8, 9)
    47 = @(7, 7)
    48 = @(7, 7)
    49 = @(7, 7)
    50 = @(7, 7)
    51 = @(7, 8)
    52 = @(8, 8)
    53 = @(8, 8)
    54 = @(5, 6)
    55 = @(6, 7)
    56 = @(8, 9)
    57 = @(7, 7)
    58 = @(8, 9)
    59 = @(6, 7)
    60 = @(6, 7)
    61 = @(7, 7)
    62 = @(6, 7)
    63 = @(7, 8)
    64 = @(6, 7)
    65 = @(4, 6)
    66 = @(3, 4)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}

